# Auto-generated edit script applying the Twintania_Profits market-data refresh
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 2179598.2
$ws.Range("J17").Value2 = 2235323.8
$ws.Range("L17").Value2 = 6705971.399999999
$ws.Range("N17").Value2 = -6706307.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 8269.625
$ws.Range("I45").Value2 = 8269.625
$ws.Range("K45").Value2 = 8269.625
$ws.Range("M45").Value2 = -7892.625
$ws.Range("H97").Value2 = 1138.0646
$ws.Range("I97").Value2 = 1075.037
$ws.Range("K97").Value2 = 1075.037
$ws.Range("M97").Value2 = -579.037
$ws.Range("H122").Value2 = 3210.2942
$ws.Range("I122").Value2 = 3121.3076
$ws.Range("K122").Value2 = 9363.9228
$ws.Range("M122").Value2 = -6913.9228
$ws.Range("H132").Value2 = 3566.5557
$ws.Range("I132").Value2 = 2071.9656
$ws.Range("J132").Value2 = 9758.429
$ws.Range("K132").Value2 = 6215.8968
$ws.Range("L132").Value2 = 29275.287
$ws.Range("M132").Value2 = -3685.8968
$ws.Range("N132").Value2 = -34335.287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 10938.023
$ws.Range("I20").Value2 = 13861.223
$ws.Range("J20").Value2 = 5676.2666
$ws.Range("K20").Value2 = 13861.223
$ws.Range("L20").Value2 = 5676.2666
$ws.Range("M20").Value2 = -13614.223
$ws.Range("N20").Value2 = -6170.2666
$ws.Range("H80").Value2 = 268.2143
$ws.Range("J80").Value2 = 300.57144
$ws.Range("L80").Value2 = 300.57144
$ws.Range("N80").Value2 = -2296.57144
$ws.Range("H83").Value2 = 268.2143
$ws.Range("J83").Value2 = 300.57144
$ws.Range("L83").Value2 = 1502.8572
$ws.Range("N83").Value2 = -11486.8572
$ws.Range("H86").Value2 = 273500
$ws.Range("J86").Value2 = 3523.5667
$ws.Range("L86").Value2 = 3523.5667
$ws.Range("N86").Value2 = -5769.566699999999
$ws.Range("H89").Value2 = 273500
$ws.Range("J89").Value2 = 3523.5667
$ws.Range("L89").Value2 = 17617.8335
$ws.Range("N89").Value2 = -28849.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value2 = 37451.555
$ws.Range("J41").Value2 = 41133
$ws.Range("L41").Value2 = 41133
$ws.Range("N41").Value2 = -41989
$ws.Range("H47").Value2 = 0
$ws.Range("I47").Value2 = 0
$ws.Range("J47").Value2 = 0
$ws.Range("K47").Value2 = 0
$ws.Range("L47").Value2 = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value2 = 500954.97
$ws.Range("I11").Value2 = 1105.3684
$ws.Range("K11").Value2 = 3316.1052
$ws.Range("M11").Value2 = -3176.1052
$ws.Range("H14").Value2 = 2916.875
$ws.Range("I14").Value2 = 2916.875
$ws.Range("K14").Value2 = 8750.625
$ws.Range("M14").Value2 = -8577.625
$ws.Range("H68").Value2 = 2835
$ws.Range("J68").Value2 = 4003
$ws.Range("L68").Value2 = 12009
$ws.Range("N68").Value2 = -13631
$ws.Range("H71").Value2 = 2835
$ws.Range("J71").Value2 = 4003
$ws.Range("L71").Value2 = 36027
$ws.Range("N71").Value2 = -44139
$ws.Range("H104").Value2 = 11944.378
$ws.Range("I104").Value2 = 7249.5
$ws.Range("K104").Value2 = 21748.5
$ws.Range("M104").Value2 = -19127.5
$ws.Range("H122").Value2 = 6667404.5
$ws.Range("I122").Value2 = 671.55554
$ws.Range("J122").Value2 = 16667504
$ws.Range("K122").Value2 = 6043.99986
$ws.Range("L122").Value2 = 150007536
$ws.Range("M122").Value2 = -3593.99986
$ws.Range("N122").Value2 = -150012436
$ws.Range("H131").Value2 = 2893.7097
$ws.Range("J131").Value2 = 3967.6843
$ws.Range("L131").Value2 = 11903.0529
$ws.Range("N131").Value2 = -21983.0529
$ws.Range("H132").Value2 = 4547188
$ws.Range("I132").Value2 = 1741
$ws.Range("J132").Value2 = 10001724
$ws.Range("K132").Value2 = 15669
$ws.Range("L132").Value2 = 90015516
$ws.Range("M132").Value2 = -13139
$ws.Range("N132").Value2 = -90020576
$ws.Range("H133").Value2 = 6125
$ws.Range("I133").Value2 = 5000
$ws.Range("K133").Value2 = 15000
$ws.Range("M133").Value2 = -9940
$ws.Range("H139").Value2 = 3379.1
$ws.Range("I139").Value2 = 3379.1
$ws.Range("J139").Value2 = 0
$ws.Range("K139").Value2 = 10137.3
$ws.Range("L139").Value2 = 0
$ws.Range("M139").Value2 = -4997.299999999999
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value2 = 0
$ws.Range("J15").Value2 = 0
$ws.Range("L15").Value2 = 0
$ws.Range("N15").ClearContents()
$ws.Range("H70").Value2 = 8434.5625
$ws.Range("I70").Value2 = 7587.364
$ws.Range("J70").Value2 = 10298.4
$ws.Range("K70").Value2 = 7587.364
$ws.Range("L70").Value2 = 10298.4
$ws.Range("M70").Value2 = -7317.364
$ws.Range("N70").Value2 = -10838.4
$ws.Range("H73").Value2 = 8434.5625
$ws.Range("I73").Value2 = 7587.364
$ws.Range("J73").Value2 = 10298.4
$ws.Range("K73").Value2 = 7587.364
$ws.Range("L73").Value2 = 10298.4
$ws.Range("M73").Value2 = -6651.364
$ws.Range("N73").Value2 = -12170.4
$ws.Range("H80").Value2 = 2477.6
$ws.Range("I80").Value2 = 2396.3333
$ws.Range("J80").Value2 = 2599.5
$ws.Range("K80").Value2 = 2396.3333
$ws.Range("L80").Value2 = 2599.5
$ws.Range("M80").Value2 = -1398.3333
$ws.Range("N80").Value2 = -4595.5
$ws.Range("H81").Value2 = 0
$ws.Range("J81").Value2 = 0
$ws.Range("L81").Value2 = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value2 = 2477.6
$ws.Range("I83").Value2 = 2396.3333
$ws.Range("J83").Value2 = 2599.5
$ws.Range("K83").Value2 = 11981.6665
$ws.Range("L83").Value2 = 12997.5
$ws.Range("M83").Value2 = -6989.666499999999
$ws.Range("N83").Value2 = -22981.5
$ws.Range("H84").Value2 = 0
$ws.Range("J84").Value2 = 0
$ws.Range("L84").Value2 = 0
$ws.Range("N84").ClearContents()
$ws.Range("H102").Value2 = 1871.1818
$ws.Range("I102").Value2 = 1871.1818
$ws.Range("K102").Value2 = 1871.1818
$ws.Range("M102").Value2 = -249.1818000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 5881.4546
$ws.Range("I40").Value2 = 5677.3335
$ws.Range("K40").Value2 = 5677.3335
$ws.Range("M40").Value2 = -5541.3335
$ws.Range("H68").Value2 = 3451.4614
$ws.Range("J68").Value2 = 3479.375
$ws.Range("L68").Value2 = 3479.375
$ws.Range("N68").Value2 = -4977.375
$ws.Range("H71").Value2 = 3451.4614
$ws.Range("J71").Value2 = 3479.375
$ws.Range("L71").Value2 = 17396.875
$ws.Range("N71").Value2 = -24884.875
$ws.Range("H122").Value2 = 3062.0417
$ws.Range("I122").Value2 = 2524.2666
$ws.Range("J122").Value2 = 3958.3333
$ws.Range("K122").Value2 = 7572.7998
$ws.Range("L122").Value2 = 11874.9999
$ws.Range("M122").Value2 = -5122.7998
$ws.Range("N122").Value2 = -16774.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value2 = 64444.332
$ws.Range("I92").Value2 = 0
$ws.Range("J92").Value2 = 64444.332
$ws.Range("K92").Value2 = 0
$ws.Range("L92").Value2 = 64444.332
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value2 = -69436.33199999999
$ws.Range("H122").Value2 = 560645.6
$ws.Range("I122").Value2 = 656289.1
$ws.Range("K122").Value2 = 1968867.3
$ws.Range("M122").Value2 = -1966417.3
$ws.Range("H126").Value2 = 8358.324000000001
$ws.Range("I126").Value2 = 5774.077
$ws.Range("K126").Value2 = 17322.231
$ws.Range("M126").Value2 = -14852.231
